# [AA | 13/4/2018] : commit for Sample tab changes and profile page
#
# Adds four new worksheets (Categorical_Data, Profile_Data, Sample_Data,
# Fields_Data) ahead of the original Sheet1, and populates them with the
# profiling output data. The original Sheet1 data is preserved verbatim
# on the new "Sample_Data" tab; the original "Sheet1" tab is left blank.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------
# The original Sheet1 held the raw sample rows; that data is reproduced
# on the new "Sample_Data" tab below, so clear it here (Sheet1 ends up
# blank in the edited workbook).
# ---------------------------------------------------------------------
$ws.Cells.Clear()

$sampleData = @(
  @(1, 100, "User",   26, "Goa",    10000),
  @(2, 200, "User2",  27, "Pune",   20000),
  @(3, 500, "User4",  27, "Mumbai", 40000),
  @(4, 222, "User54", 22, "US",     60000),
  @(5, 333, "User21", 21, "Rome",   70000),
  @(6, 444, "User42", 24, "Venice", 90000),
  @(7, 555, "User34", 25, "Hyd",    120000),
  @(8, 666, "User25", 27, "Goa",    230000),
  @(9, 777, "User48", 24, "Pune",   120000)
)

# ---------------------------------------------------------------------
# Sheets are inserted ahead of the active sheet, so adding them in
# reverse of the desired final order leaves the tabs reading left to
# right as: Categorical_Data, Profile_Data, Sample_Data, Fields_Data,
# Sheet1.
# ---------------------------------------------------------------------

# --- Fields_Data -------------------------------------------------------
$wsFields = $wb.Worksheets.Add()
$wsFields.Name = "Fields_Data"

$fieldsData = @(
  @(1, "ID",      "INTEGER"),
  @(2, "NAME",    "STRING"),
  @(3, "AGE",     "INTEGER"),
  @(4, "ADDRESS", "STRING"),
  @(5, "SALARY",  "DOUBLE")
)
for ($r = 0; $r -lt $fieldsData.Length; $r++) {
  $row = $fieldsData[$r]
  $wsFields.Cells.Item($r + 1, 1).Value = $row[0]
  $wsFields.Cells.Item($r + 1, 2).Value = $row[1]
  $wsFields.Cells.Item($r + 1, 3).Value = $row[2]
  $wsFields.Cells.Item($r + 1, 6).Value = "No"
  $wsFields.Cells.Item($r + 1, 7).Value = "No"
}

# --- Sample_Data ---------------------------------------------------------
$wsSample = $wb.Worksheets.Add()
$wsSample.Name = "Sample_Data"

for ($r = 0; $r -lt $sampleData.Length; $r++) {
  $row = $sampleData[$r]
  $wsSample.Cells.Item($r + 1, 1).Value = $row[0]
  $wsSample.Cells.Item($r + 1, 2).Value = $row[1]
  $wsSample.Cells.Item($r + 1, 3).Value = $row[2]
  $wsSample.Cells.Item($r + 1, 4).Value = $row[3]
  $wsSample.Cells.Item($r + 1, 5).Value = $row[4]
  $wsSample.Cells.Item($r + 1, 6).Value = $row[5]
}
$wsSample.Cells.Item(10, 1).Value = 10

# --- Profile_Data --------------------------------------------------------
$wsProfile = $wb.Worksheets.Add()
$wsProfile.Name = "Profile_Data"

$profileHeader = @("Field Name", "Data Type", "Scale Type", "Distinct Count", "Unique Count(%)", "Missing(%)")
for ($c = 0; $c -lt $profileHeader.Length; $c++) {
  $wsProfile.Cells.Item(1, $c + 1).Value = $profileHeader[$c]
}

$profileData = @(
  @("id",            "INT",    "categorical", 9,  90,  10),
  @("name",          "STRING", "categorical", 10, 100, 0),
  @("age",           "INT",    "categorical", 6,  60,  10),
  @("address",       "STRING", "categorical", 8,  80,  0),
  @("salary",        "DOUBLE", "categorical", 8,  80,  10),
  @("jobinstanceid", "STRING", "constant",    1,  10,  0)
)
for ($r = 0; $r -lt $profileData.Length; $r++) {
  $row = $profileData[$r]
  for ($c = 0; $c -lt $row.Length; $c++) {
    $wsProfile.Cells.Item($r + 2, $c + 1).Value = $row[$c]
  }
}

# --- Categorical_Data ------------------------------------------------------
$wsCategorical = $wb.Worksheets.Add()
$wsCategorical.Name = "Categorical_Data"

$wsCategorical.Cells.Item(1, 1).Value = "Data"
$wsCategorical.Cells.Item(1, 2).Value = "Frequency"
$wsCategorical.Cells.Item(1, 3).Value = "Percentage"

$categoricalData = @(777, 666, 555, 500, 444, 333, 222, 200, 100)
for ($r = 0; $r -lt $categoricalData.Length; $r++) {
  $rowIdx = $r + 2
  $wsCategorical.Cells.Item($rowIdx, 1).Value = $categoricalData[$r]
  $wsCategorical.Cells.Item($rowIdx, 2).Value = 1
  $wsCategorical.Cells.Item($rowIdx, 3).Value = 0.1111
  $wsCategorical.Cells.Item($rowIdx, 3).NumberFormat = "0.00%"
}

# Leave the original tab ("Sheet1") blank, and make the newly-created
# first tab the active/selected one, matching the authored workbook.
$wsCategorical.Activate()
